$wb = $excel.ActiveWorkbook

$dateFormat = "yyyy-mm-dd HH:mm:ss"

# Overview sheet: G2 holds the "Latest HO Xliff Generate Date" for the first file.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-27 19:12:14"
$wsOverview.Range("G2").NumberFormat = $dateFormat

# zh-cn sheet: H2 = Correspond Handoff Datetime, K2 = Correspond Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-27 19:12:09"
$wsZhCn.Range("H2").NumberFormat = $dateFormat
$wsZhCn.Range("K2").Value = "2016-08-27 19:12:26"
$wsZhCn.Range("K2").NumberFormat = $dateFormat

# de-de sheet: H2 = Correspond Handoff Datetime, K2 = Correspond Handback DateTime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-27 19:12:14"
$wsDeDe.Range("H2").NumberFormat = $dateFormat
$wsDeDe.Range("K2").Value = "2016-08-27 19:12:33"
$wsDeDe.Range("K2").NumberFormat = $dateFormat
